# DALA-3260: update lksg – realign "Component" column (F) shared-string
# references for the EU Taxonomy "AmountWithCurrency" rows and move the
# selection/scroll position, mirroring the authored xlsx diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Component cell (column F) held
# "Custom EuTaxonomyAmountWithCurrencyComponent" and are renamed to the
# shorter "AmountWithCurrencyComponent" label.
$rowsToRename = @(12, 14, 16, 19, 31, 33, 35, 38, 50, 52, 54, 57)
foreach ($r in $rowsToRename) {
    $ws.Cells.Item($r, 6).Value = "AmountWithCurrencyComponent"
}

# Update the view state: scroll position and active selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("F57").Select()
